$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3050.4
$ws.Range("I43").Value = 2416.3333
$ws.Range("K43").Value = 2416.3333
$ws.Range("M43").Value = -2347.3333
$ws.Range("H62").Value = 3391.1667
$ws.Range("I62").Value = 3186.889
$ws.Range("K62").Value = 3186.889
$ws.Range("M62").Value = -2562.889
$ws.Range("H64").Value = 26323642
$ws.Range("I64").Value = 7981.032
$ws.Range("J64").Value = 142864430
$ws.Range("K64").Value = 7981.032
$ws.Range("L64").Value = 142864430
$ws.Range("M64").Value = -7733.032
$ws.Range("N64").Value = -142864926
$ws.Range("H65").Value = 3391.1667
$ws.Range("I65").Value = 3186.889
$ws.Range("K65").Value = 15934.445
$ws.Range("M65").Value = -12814.445
$ws.Range("H67").Value = 26323642
$ws.Range("I67").Value = 7981.032
$ws.Range("J67").Value = 142864430
$ws.Range("K67").Value = 7981.032
$ws.Range("L67").Value = 142864430
$ws.Range("M67").Value = -7123.032
$ws.Range("N67").Value = -142866146
$ws.Range("H116").Value = 6590.5
$ws.Range("I116").Value = 6309.2
$ws.Range("K116").Value = 6309.2
$ws.Range("M116").Value = -2867.2
$ws.Range("H132").Value = 3530.7222
$ws.Range("I132").Value = 3672.1042
$ws.Range("K132").Value = 11016.3126
$ws.Range("M132").Value = -8486.312600000001
$ws.Range("H135").Value = 2575.125
$ws.Range("I135").Value = 2560.1333
$ws.Range("K135").Value = 23041.1997
$ws.Range("M135").Value = -20506.1997
$ws.Range("H137").Value = 964468.0600000001
$ws.Range("I137").Value = 2501721.5
$ws.Range("K137").Value = 7505164.5
$ws.Range("M137").Value = -7502614.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6041.5557
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 6041.5557
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6041.5557
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -6615.5557
$ws.Range("H45").Value = 32533.285
$ws.Range("I45").Value = 43141
$ws.Range("K45").Value = 43141
$ws.Range("M45").Value = -42764
$ws.Range("H61").Value = 1992.1305
$ws.Range("I61").Value = 686.0769
$ws.Range("K61").Value = 686.0769
$ws.Range("M61").Value = -474.0769
$ws.Range("H74").Value = 193625.28
$ws.Range("I74").Value = 428116.94
$ws.Range("K74").Value = 428116.94
$ws.Range("M74").Value = -427242.94
$ws.Range("H77").Value = 193625.28
$ws.Range("I77").Value = 428116.94
$ws.Range("K77").Value = 2140584.7
$ws.Range("M77").Value = -2136216.7
$ws.Range("H110").Value = 3328.3076
$ws.Range("I110").Value = 982.1667
$ws.Range("J110").Value = 5339.2856
$ws.Range("K110").Value = 982.1667
$ws.Range("L110").Value = 5339.2856
$ws.Range("M110").Value = 1062.8333
$ws.Range("N110").Value = -9429.285599999999
$ws.Range("H122").Value = 8798.625
$ws.Range("I122").Value = 5065.6665
$ws.Range("K122").Value = 15196.9995
$ws.Range("M122").Value = -12746.9995
$ws.Range("H136").Value = 1992.1305
$ws.Range("I136").Value = 686.0769
$ws.Range("K136").Value = 2058.2307
$ws.Range("M136").Value = 491.7692999999999
$ws.Range("H139").Value = 76725.45
$ws.Range("J139").Value = 76725.45
$ws.Range("L139").Value = 76725.45
$ws.Range("N139").Value = -87005.45
$ws.Range("H141").Value = 98333.336
$ws.Range("J141").Value = 98333.336
$ws.Range("L141").Value = 98333.336
$ws.Range("N141").Value = -108693.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 735
$ws.Range("I64").Value = 330.33334
$ws.Range("K64").Value = 330.33334
$ws.Range("M64").Value = -105.33334
$ws.Range("H67").Value = 735
$ws.Range("I67").Value = 330.33334
$ws.Range("K67").Value = 330.33334
$ws.Range("M67").Value = 449.66666
$ws.Range("H86").Value = 2321.389
$ws.Range("I86").Value = 2209.375
$ws.Range("K86").Value = 2209.375
$ws.Range("M86").Value = -1086.375
$ws.Range("H89").Value = 2321.389
$ws.Range("I89").Value = 2209.375
$ws.Range("K89").Value = 11046.875
$ws.Range("M89").Value = -5430.875
$ws.Range("H94").Value = 74075260
$ws.Range("I94").Value = 100000370
$ws.Range("K94").Value = 100000370
$ws.Range("M94").Value = -99999919
$ws.Range("H97").Value = 10434
$ws.Range("I97").Value = 10434
$ws.Range("K97").Value = 10434
$ws.Range("M97").Value = -9443
$ws.Range("H134").Value = 4328.6177
$ws.Range("I134").Value = 4435.136
$ws.Range("J134").Value = 4133.3335
$ws.Range("K134").Value = 13305.408
$ws.Range("L134").Value = 12400.0005
$ws.Range("M134").Value = -10770.408
$ws.Range("N134").Value = -17470.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1807.5454
$ws.Range("I22").Value = 1788.3
$ws.Range("K22").Value = 1788.3
$ws.Range("M22").Value = -1438.3
$ws.Range("H31").Value = 4849.577
$ws.Range("I31").Value = 3137.647
$ws.Range("K31").Value = 3137.647
$ws.Range("M31").Value = -2842.647
$ws.Range("H34").Value = 4849.577
$ws.Range("I34").Value = 3137.647
$ws.Range("K34").Value = 3137.647
$ws.Range("M34").Value = -2935.647
$ws.Range("H88").Value = 42875
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 42875
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 42875
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -43687
$ws.Range("H91").Value = 42875
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 42875
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 42875
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -45683
$ws.Range("H105").Value = 1118.5
$ws.Range("I105").Value = 725.7857
$ws.Range("K105").Value = 725.7857
$ws.Range("M105").Value = 1021.2143
$ws.Range("H107").Value = 965.44446
$ws.Range("I107").Value = 407.5
$ws.Range("K107").Value = 407.5
$ws.Range("M107").Value = 1512.5
$ws.Range("H132").Value = 1829.3158
$ws.Range("I132").Value = 1517.2667
$ws.Range("K132").Value = 4551.800099999999
$ws.Range("M132").Value = -2021.800099999999
$ws.Range("H134").Value = 1931.84
$ws.Range("I134").Value = 1702.0444
$ws.Range("K134").Value = 5106.1332
$ws.Range("M134").Value = -2571.1332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1118.2222
$ws.Range("I5").Value = 701.7
$ws.Range("J5").Value = 1638.875
$ws.Range("K5").Value = 2105.1
$ws.Range("L5").Value = 4916.625
$ws.Range("M5").Value = -1993.1
$ws.Range("N5").Value = -5140.625
$ws.Range("H8").Value = 737.5
$ws.Range("I8").Value = 737.5
$ws.Range("K8").Value = 2212.5
$ws.Range("M8").Value = -2073.5
$ws.Range("H107").Value = 790.6
$ws.Range("J107").Value = 584.3333
$ws.Range("L107").Value = 1752.9999
$ws.Range("N107").Value = -5592.9999
$ws.Range("H131").Value = 4401.2744
$ws.Range("J131").Value = 1675.8292
$ws.Range("L131").Value = 5027.487599999999
$ws.Range("N131").Value = -15107.4876
$ws.Range("H135").Value = 1118.2222
$ws.Range("I135").Value = 701.7
$ws.Range("J135").Value = 1638.875
$ws.Range("K135").Value = 6315.3
$ws.Range("L135").Value = 14749.875
$ws.Range("M135").Value = -3780.3
$ws.Range("N135").Value = -19819.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 106023.5
$ws.Range("I70").Value = 226065.22
$ws.Range("K70").Value = 226065.22
$ws.Range("M70").Value = -225795.22
$ws.Range("H73").Value = 106023.5
$ws.Range("I73").Value = 226065.22
$ws.Range("K73").Value = 226065.22
$ws.Range("M73").Value = -225129.22
$ws.Range("H107").Value = 993
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 3763.5588
$ws.Range("I122").Value = 2543.5881
$ws.Range("J122").Value = 4983.5293
$ws.Range("K122").Value = 7630.7643
$ws.Range("L122").Value = 14950.5879
$ws.Range("M122").Value = -5180.7643
$ws.Range("N122").Value = -19850.5879
$ws.Range("H132").Value = 1230.1052
$ws.Range("I132").Value = 898.25
$ws.Range("K132").Value = 2694.75
$ws.Range("M132").Value = -164.75
$ws.Range("H136").Value = 10347.667
$ws.Range("J136").Value = 10347.667
$ws.Range("L136").Value = 31043.001
$ws.Range("N136").Value = -36143.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 756.53845
$ws.Range("I93").Value = 756.53845
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 756.53845
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 491.46155
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 942.6
$ws.Range("I107").Value = 969.55554
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 2908.66662
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = -988.66662
$ws.Range("N107").Value = -5940
$ws.Range("H126").Value = 2777.7856
$ws.Range("I126").Value = 1859
$ws.Range("K126").Value = 5577
$ws.Range("M126").Value = -3107
$ws.Range("H132").Value = 2922.3333
$ws.Range("I132").Value = 2808.8948
$ws.Range("K132").Value = 8426.6844
$ws.Range("M132").Value = -5896.6844
$ws.Range("H136").Value = 1056.375
$ws.Range("I136").Value = 908.3
$ws.Range("J136").Value = 1303.1666
$ws.Range("K136").Value = 2724.9
$ws.Range("L136").Value = 3909.4998
$ws.Range("M136").Value = -174.8999999999996
$ws.Range("N136").Value = -9009.4998
